# Do Log In to ACME
# Adds a new "Framework\KillAllProcesses.xaml" test case to the Tests sheet
# and populates the Result sheet with the (copied) test outcomes.

$wb = $excel.ActiveWorkbook
$wsTests = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

# --- Add the new test case row to the "Tests" sheet ---
$wsTests.Cells.Item(10, 1).Value = "Framework\KillAllProcesses.xaml"
$wsTests.Cells.Item(10, 2).Value = "Success"

# --- Populate the "Result" sheet with the test run results ---
$results = @(
    @("Framework\InitAllSettings.xaml", "Success"),
    @("Framework\InitAllApplications.xaml", "Success"),
    @("Framework\CloseAllApplications.xaml", "Success"),
    @("Framework\CloseAllApplications.xaml", "SystemException"),
    @("Framework\InitAllSettings.xaml", "Success"),
    @("Framework\InitAllSettings.xaml", "Success"),
    @("Framework\InitAllApplications.xaml", "Success"),
    @("Framework\CloseAllApplications.xaml", "Success"),
    @("Framework\KillAllProcesses.xaml", "Success")
)

$row = 2
foreach ($entry in $results) {
    $wsResult.Cells.Item($row, 1).Value = $entry[0]
    $wsResult.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

# --- Update selections / active sheet to match the recorded view state ---
$wsTests.Range("A30").Select()
$wsResult.Activate()
$wsResult.Range("D18").Select()
